$wb = $excel.ActiveWorkbook

# --- zh-cn sheet: row 3 (5b2edd21-ede1-464f-9261-522a3c1ef648 file) is now handed back ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("K3").Value = "2016-09-06 08:16:43"
$wsZh.Range("P3").Value = ""
$wsZh.Columns.Item(16).AutoFit()

# --- de-de sheet: row 3 (5b2edd21-ede1-464f-9261-522a3c1ef648 file) is now handed back ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("K3").Value = "2016-09-06 08:16:59"
$wsDe.Range("P3").Value = ""
$wsDe.Columns.Item(16).AutoFit()

# --- Overview sheet: row 3 reflects the new combined status for both locales ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
